$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the row for even_MAG-GUT66322.fa entirely (old row 3), shifting rows up.
$ws.Rows.Item(3).Delete()

# Delete the "max" column (old column C), shifting "prediction"/"rejection-f" left.
$ws.Columns.Item(3).Delete()

# Update the remaining data rows (now rows 2-5) with the new values.
$ws.Range("A2").Value = "even_MAG-GUT43773.fa"
$ws.Range("B2").Value = -6.879024153281895
$ws.Range("C2").Value = "s__CAG-288 sp000437395"
$ws.Range("D2").Value = "s__CAG-288 sp000437395(reject)"

$ws.Range("A3").Value = "even_MAG-GUT66408.fa"
$ws.Range("B3").Value = -3.034966611753012
$ws.Range("C3").Value = "s__CAG-288 sp000437395"
$ws.Range("D3").Value = "s__CAG-288 sp000437395(reject)"

$ws.Range("A4").Value = "even_MAG-GUT70664.fa"
$ws.Range("B4").Value = -1.822538634812249
$ws.Range("C4").Value = "s__CAG-288 sp000437395"
$ws.Range("D4").Value = "s__CAG-288 sp000437395(reject)"

$ws.Range("A5").Value = "even_MAG-GUT71751.fa"
$ws.Range("B5").Value = -1.218436119744123
$ws.Range("C5").Value = "s__CAG-288 sp000437395"
$ws.Range("D5").Value = "s__CAG-288 sp000437395(reject)"
